$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestDataSheet1")

# Change C4 from the text "ssf" to the number 12
$ws.Range("C4").Value = 12

# Move the active selection to D4 (was N19)
$ws.Range("D4").Select()
